$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.395.21"
$ws.Range("E2").Value = "  -5.05%  "

# Row 3
$ws.Range("D3").Value = "3.303.09"
$ws.Range("E3").Value = "  -6.92%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "'181.98"
$ws.Range("E5").Value = "  -11.17%  "

# Row 6
$ws.Range("D6").Value = "'530.14"
$ws.Range("E6").Value = "  -4.73%  "

# Row 7
$ws.Range("D7").Value = "'0.605"
$ws.Range("E7").Value = "  -0.55%  "

# Row 8
$ws.Range("D8").Value = "3.294.89"
$ws.Range("E8").Value = "  -6.86%  "

# Row 9
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").Value = "'0.622"
$ws.Range("E10").Value = "  -6.21%  "

# Row 11
$ws.Range("D11").Value = "'59.15"
$ws.Range("E11").Value = "  -8.35%  "

# Row 12
$ws.Range("E12").Value = "  -7.29%  "

# Row 13
$ws.Range("D13").Value = "'0.0000263"
$ws.Range("E13").Value = "  -3.13%  "

# Row 14
$ws.Range("D14").Value = "'9.14"
$ws.Range("E14").Value = "  -8.00%  "

# Row 15
$ws.Range("D15").Value = "3.824.95"
$ws.Range("E15").Value = "  -7.25%  "

# Row 16
$ws.Range("D16").Value = "3.300.55"
$ws.Range("E16").Value = "  -7.10%  "

# Row 17
$ws.Range("E17").Value = "  -5.71%  "

# Row 18
$ws.Range("D18").Value = "'17.77"
$ws.Range("E18").Value = "  -5.01%  "

# Row 19
$ws.Range("D19").Value = "64.174.81"
$ws.Range("E19").Value = "  -5.19%  "

# Row 20
$ws.Range("D20").Value = "'11.13"
$ws.Range("E20").Value = "  -7.04%  "

# Row 21
$ws.Range("D21").Value = "'0.964"
$ws.Range("E21").Value = "  -7.58%  "

# Row 22
$ws.Range("D22").Value = "'374.39"
$ws.Range("E22").Value = "  -5.59%  "

# Row 23
$ws.Range("D23").Value = "'3.83"
$ws.Range("E23").Value = "  -6.05%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'81.21"
$ws.Range("E24").Value = "  -2.36%  "

# Row 25
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "'11.22"
$ws.Range("E25").Value = "  -7.19%  "

# Row 26
$ws.Range("D26").Value = "'3.93"
$ws.Range("E26").Value = "  +3.97%  "

# Row 27
$ws.Range("E27").Value = "  -1.34%  "

# Row 28
$ws.Range("D28").Value = "'2.68"
$ws.Range("E28").Value = "  -5.31%  "

# Row 29
$ws.Range("D29").Value = "'11.64"
$ws.Range("E29").Value = "  -5.19%  "

# Row 30
$ws.Range("D30").Value = "'8.45"
$ws.Range("E30").Value = "  -5.86%  "

# Row 31
$ws.Range("D31").Value = "'29.04"
$ws.Range("E31").Value = "  -6.72%  "

# Row 32
$ws.Range("D32").Value = "'6.82"
$ws.Range("E32").Value = "  -6.88%  "

# Row 33
$ws.Range("D33").Value = "'643.91"

# Row 34
$ws.Range("D34").Value = "'11.35"
$ws.Range("E34").Value = "  -4.63%  "

# Row 35
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'59.53"
$ws.Range("E35").Value = "  -7.39%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.105"
$ws.Range("E36").Value = "  -6.25%  "

# Row 37
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.396"
$ws.Range("E38").Value = "  -4.16%  "

# Row 39
$ws.Range("D39").Value = "'36.89"
$ws.Range("E39").Value = "  -6.87%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0728"
$ws.Range("E40").Value = "  +3.44%  "

# Row 41
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  -0.16%  "

# Row 42
$ws.Range("D42").Value = "'0.128"
$ws.Range("E42").Value = "  -3.52%  "

# Row 43
$ws.Range("D43").Value = "2.906.68"
$ws.Range("E43").Value = "  -6.48%  "

# Row 44
$ws.Range("E44").Value = "  -2.84%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0407"
$ws.Range("E45").Value = "  +0.56%  "

# Row 46
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.71"
$ws.Range("E46").Value = "  -11.87%  "

# Row 47
$ws.Range("E47").Value = "  +9.55%  "

# Row 48
$ws.Range("D48").Value = "'2.67"
$ws.Range("E48").Value = "  -2.37%  "

# Row 49
$ws.Range("D49").Value = "'2.61"
$ws.Range("E49").Value = "  -9.91%  "

# Row 50
$ws.Range("E50").Value = "  -1.07%  "

# Row 51
$ws.Range("E51").Value = "  -0.93%  "
